# Updated cryptos list on Sun Apr 30 04:34:16 UTC 2023 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for the cryptos sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Force the literal text into the cell (not auto-coerced to a number),
    # then drop back to the default style so no stray formatting is left behind.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "29.522.02"
Set-TextValue $ws.Range("E2") "  -0.05%  "
Set-TextValue $ws.Range("D3") "1.918.24"
Set-TextValue $ws.Range("E3") "  +0.35%  "
Set-TextValue $ws.Range("E4") "  +0.85%  "
Set-TextValue $ws.Range("D5") "325.51"
Set-TextValue $ws.Range("E5") "  -0.27%  "
Set-TextValue $ws.Range("E6") "  +0.75%  "
Set-TextValue $ws.Range("D7") "0.4812"
Set-TextValue $ws.Range("E7") "  -0.67%  "
Set-TextValue $ws.Range("D8") "0.4052"
Set-TextValue $ws.Range("E8") "  -0.68%  "
Set-TextValue $ws.Range("D9") "0.08215"
Set-TextValue $ws.Range("E9") "  +0.78%  "
Set-TextValue $ws.Range("E10") "  -0.49%  "
Set-TextValue $ws.Range("D11") "23.41"
Set-TextValue $ws.Range("E11") "  -0.44%  "
Set-TextValue $ws.Range("D12") "1.926.42"
Set-TextValue $ws.Range("E12") "  +0.97%  "
Set-TextValue $ws.Range("D13") "6.052"
Set-TextValue $ws.Range("E13") "  +0.35%  "
Set-TextValue $ws.Range("D14") "7.236"
Set-TextValue $ws.Range("E14") "  +1.74%  "
Set-TextValue $ws.Range("D15") "91.43"
Set-TextValue $ws.Range("E15") "  +1.12%  "
Set-TextValue $ws.Range("D16") "0.06871"
Set-TextValue $ws.Range("E16") "  +1.59%  "
Set-TextValue $ws.Range("E17") "  +0.70%  "
Set-TextValue $ws.Range("D18") "0.00001038"
Set-TextValue $ws.Range("E18") "  -0.30%  "
Set-TextValue $ws.Range("E19") "  -1.15%  "
Set-TextValue $ws.Range("E20") "  +0.65%  "
Set-TextValue $ws.Range("D21") "29.522.22"
Set-TextValue $ws.Range("E21") "  -0.08%  "
Set-TextValue $ws.Range("D23") "11.87"
Set-TextValue $ws.Range("E23") "  +0.43%  "
Set-TextValue $ws.Range("D24") "2.193"
Set-TextValue $ws.Range("E24") "  +1.27%  "
Set-TextValue $ws.Range("D25") "2.145.24"
Set-TextValue $ws.Range("E25") "  +0.30%  "
Set-TextValue $ws.Range("D26") "6.532"
Set-TextValue $ws.Range("E26") "  +3.82%  "
Set-TextValue $ws.Range("D27") "156.01"
Set-TextValue $ws.Range("E27") "  +0.78%  "
Set-TextValue $ws.Range("E28") "  -0.21%  "
Set-TextValue $ws.Range("D29") "2.095"
Set-TextValue $ws.Range("E29") "  -0.63%  "
Set-TextValue $ws.Range("D30") "120.58"
Set-TextValue $ws.Range("E30") "  +0.64%  "
Set-TextValue $ws.Range("E31") "  -1.71%  "
Set-TextValue $ws.Range("D32") "0.09637"
Set-TextValue $ws.Range("E32") "  +0.78%  "
Set-TextValue $ws.Range("D33") "5.616"
Set-TextValue $ws.Range("E33") "  +1.39%  "
Set-TextValue $ws.Range("E34") "  +0.22%  "
Set-TextValue $ws.Range("E35") "  -1.91%  "
Set-TextValue $ws.Range("D37") "0.02284"
Set-TextValue $ws.Range("E37") "  +0.61%  "
Set-TextValue $ws.Range("E38") "  +0.99%  "
Set-TextValue $ws.Range("D39") "0.5935"
Set-TextValue $ws.Range("E39") "  -0.35%  "
Set-TextValue $ws.Range("D40") "10.68"
Set-TextValue $ws.Range("E40") "  -0.01%  "
Set-TextValue $ws.Range("D41") "7.899"
Set-TextValue $ws.Range("E41") "  -0.58%  "
Set-TextValue $ws.Range("E42") "  -0.63%  "
Set-TextValue $ws.Range("D43") "2.476"
Set-TextValue $ws.Range("E43") "  +1.29%  "
Set-TextValue $ws.Range("D44") "1.246"
Set-TextValue $ws.Range("E44") "  -2.94%  "
Set-TextValue $ws.Range("D45") "12.40"
Set-TextValue $ws.Range("E45") "  -0.13%  "
Set-TextValue $ws.Range("D46") "0.07469"
Set-TextValue $ws.Range("E46") "  -3.20%  "
Set-TextValue $ws.Range("D47") "0.5559"
Set-TextValue $ws.Range("E47") "  -0.34%  "
Set-TextValue $ws.Range("D48") "1.940"
Set-TextValue $ws.Range("E48") "  -0.83%  "
Set-TextValue $ws.Range("E49") "  +2.77%  "
Set-TextValue $ws.Range("D50") "2.428"
Set-TextValue $ws.Range("E50") "  +3.27%  "
Set-TextValue $ws.Range("D51") "72.04"
Set-TextValue $ws.Range("E51") "  -1.01%  "
